$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace all occurrences of the vehicle code "VG_GT3P5" with "LOR_GT3P5"
# in column C (vehicle). This affects rows 2, 4, 6, ..., 32.
for ($r = 2; $r -le 32; $r += 2) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "VG_GT3P5") {
        $cell.Value = "LOR_GT3P5"
    }
}
